$wb = $excel.ActiveWorkbook

# Sheet "展览" (index 1): update 想去人数 values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 12
$ws1.Range("F6").Value = 5260
$ws1.Range("F10").Value = 355
$ws1.Range("F12").Value = 66

# Sheet "全部类型" (index 4): same events appear here too, mirror updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 12
$ws4.Range("F10").Value = 5260
$ws4.Range("F15").Value = 355
$ws4.Range("F17").Value = 66
